# Add 2022-Q1 data:
#  - the existing "总计" sheet (aggregate totals) is renamed to "2022-Q1" and
#    repopulated with the fund-level holdings detail for 2022-Q1.
#  - a brand-new "总计" sheet is appended after it, containing the aggregate
#    totals table (same as before) with a new top row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Repurpose the old "总计" sheet -> "2022-Q1" fund holdings detail
# ---------------------------------------------------------------------------
$detail = $wb.Worksheets.Item("总计")
$detail.Name = "2022-Q1"

# Header row (row 1): reuse the existing header formatting (style already on
# B1:D1) and extend it across to H1.
$detail.Range("D1").Copy()
$detail.Range("E1:H1").PasteSpecial(-4122)

$detail.Cells.Item(1, 2).Value = "基金代码"
$detail.Cells.Item(1, 3).Value = "基金名称"
$detail.Cells.Item(1, 4).Value = "基金规模"
$detail.Cells.Item(1, 5).Value = "股票总仓位"
$detail.Cells.Item(1, 6).Value = "仓位占比"
$detail.Cells.Item(1, 7).Value = "持有市值(亿元)"
$detail.Cells.Item(1, 8).Value = "仓位排名"

# Text-formatted columns B..G hold fund codes / names / numbers-as-text in
# the source data (codes such as "009556" must keep their leading zero, so
# they cannot be left to auto-detect as numbers).
$detail.Range("B2:G19").NumberFormat = "@"

# Extend the existing "row index" style (already on A2:A4) down to A19.
$detail.Range("A2").Copy()
$detail.Range("A2:A19").PasteSpecial(-4122)

$rows = @(
    @(0,  "900010", "中信卓越成长两年持有期混合A",             "133.02", "93.07", "2.99", "3.9773", 8),
    @(1,  "009556", "兴全合丰三年持有期混合",                   "84.68",  "92.95", "4.31", "3.6497", 8),
    @(2,  "900090", "中信卓越成长两年持有期混合B",             "86.95",  "93.07", "2.99", "2.5998", 8),
    @(3,  "004424", "汇添富文体娱乐主题混合",                   "18.52",  "90.11", "5.02", "0.9297", 3),
    @(4,  "010326", "博时消费创新混合A",                       "19.00",  "82.83", "3.79", "0.7201", 8),
    @(5,  "005644", "广发沪港深行业龙头混合",                   "13.85",  "88.27", "4.51", "0.6246", 7),
    @(6,  "001764", "广发沪港深新机遇股票",                     "11.12",  "92.18", "5.25", "0.5838", 7),
    @(7,  "006595", "广发港股通优质增长混合",                   "8.53",   "86.63", "5.09", "0.4342", 7),
    @(8,  "900100", "中信卓越成长两年持有期混合C",             "6.91",   "93.07", "2.99", "0.2066", 8),
    @(9,  "011574", "鹏华领航一年持有期混合型证券投资基金A",   "2.75",   "90.27", "6.40", "0.1760", 3),
    @(10, "008134", "鹏华优选价值股票",                         "1.86",   "92.62", "7.52", "0.1399", 3),
    @(11, "006671", "广发消费升级股票",                         "2.80",   "92.32", "4.91", "0.1375", 7),
    @(12, "006136", "广发估值优势混合A",                       "1.96",   "94.78", "5.30", "0.1039", 6),
    @(13, "010327", "博时消费创新混合C",                       "2.60",   "82.83", "3.79", "0.0985", 8),
    @(14, "011575", "鹏华领航一年持有期混合型证券投资基金C",   "1.22",   "90.27", "6.40", "0.0781", 3),
    @(15, "011969", "建信港股通精选混合A",                     "1.01",   "57.64", "3.60", "0.0364", 10),
    @(16, "011970", "建信港股通精选混合C",                     "0.33",   "57.64", "3.60", "0.0119", 10),
    @(17, "011430", "广发估值优势混合C",                       "0.03",   "94.78", "5.30", "0.0016", 6)
)

$r = 2
foreach ($row in $rows) {
    $detail.Cells.Item($r, 1).Value = $row[0]
    $detail.Cells.Item($r, 2).Value = $row[1]
    $detail.Cells.Item($r, 3).Value = $row[2]
    $detail.Cells.Item($r, 4).Value = $row[3]
    $detail.Cells.Item($r, 5).Value = $row[4]
    $detail.Cells.Item($r, 6).Value = $row[5]
    $detail.Cells.Item($r, 7).Value = $row[6]
    $detail.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Drop the "@" number-format override now that the text values are locked
# in - the data cells keep their text type but fall back to the workbook's
# default (unstyled) appearance, matching the other quarter sheets.
$detail.Range("B2:G19").Style = "Normal"

# ---------------------------------------------------------------------------
# 2) Brand-new "总计" sheet after "2022-Q1", holding the aggregate totals
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $detail)
$total.Name = "总计"

$detail.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$detail.Range("A2").Copy()
$total.Range("A2:A5").PasteSpecial(-4122)

$totalRows = @(
    @(0, "2022-Q1", 18, 14.51),
    @(1, "2021-Q4", 11, 10.37),
    @(2, "2021-Q3", 5,  10.86),
    @(3, "2021-Q2", 9,  6.29)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

Write-Host "done"
